# Updated cryptos list on Sun Jan 21 04:44:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.647.93"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "2.473.69"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.61%  "
$ws.Range("E9").Value = "  +2.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0867"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.03"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "2.853.71"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.63%  "
$ws.Range("D16").Value = "2.478.08"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "41.601.20"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0955"
$ws.Range("E20").Value = "  +1.71%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.72"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.76%  "
$ws.Range("E28").Value = "  +3.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.09"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0772"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("E35").Value = "  +0.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +1.45%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("E42").Value = "  +2.08%  "
$ws.Range("D43").Value = "1.985.98"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0287"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.25"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.11%  "
$ws.Range("D48").Value = "2.710.39"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.52"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.25%  "
